$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D2").Value = 722
$ws.Range("E2").Value = 1070
$ws.Range("F2").Value = 1375
$ws.Range("G2").Value = 295.5
$ws.Range("H2").Value = 9846715
$ws.Range("I2").Value = "KSCP"

$ws.Range("D3").Value = 722
$ws.Range("E3").Value = 1070
$ws.Range("F3").Value = 1375
$ws.Range("G3").Value = 295.5
$ws.Range("H3").Value = 9846715
$ws.Range("I3").Value = "KSCP"

$ws.Range("D4").Value = 722
$ws.Range("E4").Value = 1070
$ws.Range("F4").Value = 1375
$ws.Range("G4").Value = 295.5
$ws.Range("H4").Value = 9846715
$ws.Range("I4").Value = "KSCP"

$ws.Range("D5").Value = 722
$ws.Range("E5").Value = 1070
$ws.Range("F5").Value = 1375
$ws.Range("G5").Value = 295.5
$ws.Range("H5").Value = 9846715
$ws.Range("I5").Value = "KSCP"

$ws.Range("D6").Value = 722
$ws.Range("E6").Value = 1070
$ws.Range("F6").Value = 1375
$ws.Range("G6").Value = 295.5
$ws.Range("H6").Value = 9846715
$ws.Range("I6").Value = "KSCP"

$ws.Range("D7").Value = 722
$ws.Range("E7").Value = 1070
$ws.Range("F7").Value = 1375
$ws.Range("G7").Value = 295.5
$ws.Range("H7").Value = 9846715
$ws.Range("I7").Value = "KSCP"

$ws.Range("D8").Value = 722
$ws.Range("E8").Value = 1070
$ws.Range("F8").Value = 1375
$ws.Range("G8").Value = 295.5
$ws.Range("H8").Value = 9846715
$ws.Range("I8").Value = "KSCP"

$ws.Range("D9").Value = 722
$ws.Range("E9").Value = 1070
$ws.Range("F9").Value = 1375
$ws.Range("G9").Value = 295.5
$ws.Range("H9").Value = 9846715
$ws.Range("I9").Value = "KSCP"

$ws.Range("D10").Value = 722
$ws.Range("E10").Value = 1070
$ws.Range("F10").Value = 1375
$ws.Range("G10").Value = 295.5
$ws.Range("H10").Value = 9846715
$ws.Range("I10").Value = "KSCP"

$ws.Range("D11").Value = 722
$ws.Range("E11").Value = 1070
$ws.Range("F11").Value = 1375
$ws.Range("G11").Value = 295.5
$ws.Range("H11").Value = 9846715
$ws.Range("I11").Value = "KSCP"

$ws.Range("D12").Value = 722
$ws.Range("E12").Value = 1070
$ws.Range("F12").Value = 1375
$ws.Range("G12").Value = 295.5
$ws.Range("H12").Value = 9846715
$ws.Range("I12").Value = "KSCP"

$ws.Range("D13").Value = 722
$ws.Range("E13").Value = 1070
$ws.Range("F13").Value = 1375
$ws.Range("G13").Value = 295.5
$ws.Range("H13").Value = 9846715
$ws.Range("I13").Value = "KSCP"

$ws.Range("D14").Value = 722
$ws.Range("E14").Value = 1070
$ws.Range("F14").Value = 1375
$ws.Range("G14").Value = 295.5
$ws.Range("H14").Value = 9846715
$ws.Range("I14").Value = "KSCP"

$ws.Range("D15").Value = 722
$ws.Range("E15").Value = 1070
$ws.Range("F15").Value = 1375
$ws.Range("G15").Value = 295.5
$ws.Range("H15").Value = 9846715
$ws.Range("I15").Value = "KSCP"

$ws.Range("D16").Value = 722
$ws.Range("E16").Value = 1070
$ws.Range("F16").Value = 1375
$ws.Range("G16").Value = 295.5
$ws.Range("H16").Value = 9846715
$ws.Range("I16").Value = "KSCP"

$ws.Range("D17").Value = 253.5
$ws.Range("E17").Value = 200.5
$ws.Range("F17").Value = 262
$ws.Range("G17").Value = 195.5
$ws.Range("H17").Value = 9846715
$ws.Range("I17").Value = "KSCP"

$ws.Range("D18").Value = 151
$ws.Range("E18").Value = 150
$ws.Range("F18").Value = 195
$ws.Range("G18").Value = 129.5
$ws.Range("H18").Value = 9846715
$ws.Range("I18").Value = "KSCP"

$ws.Range("D19").Value = 122
$ws.Range("E19").Value = 168
$ws.Range("F19").Value = 172.5
$ws.Range("G19").Value = 115
$ws.Range("H19").Value = 9846715
$ws.Range("I19").Value = "KSCP"

$ws.Range("D20").Value = 94.5
$ws.Range("E20").Value = 79.5
$ws.Range("F20").Value = 100
$ws.Range("G20").Value = 76.75
$ws.Range("H20").Value = 9846715
$ws.Range("I20").Value = "KSCP"

$ws.Range("D21").Value = 45
$ws.Range("E21").Value = 30.70000076293945
$ws.Range("F21").Value = 45
$ws.Range("G21").Value = 27.64999961853028
$ws.Range("H21").Value = 9846715
$ws.Range("I21").Value = "KSCP"

$ws.Range("D22").Value = 27.45000076293945
$ws.Range("E22").Value = 71.5
$ws.Range("F22").Value = 112
$ws.Range("G22").Value = 26.60000038146973
$ws.Range("H22").Value = 9846715
$ws.Range("I22").Value = "KSCP"

$ws.Range("D23").Value = 40.29999923706055
$ws.Range("E23").Value = 33.09999847412109
$ws.Range("F23").Value = 49.95000076293945
$ws.Range("G23").Value = 30.54999923706055
$ws.Range("H23").Value = 9846715
$ws.Range("I23").Value = "KSCP"

$ws.Range("D24").Value = 30
$ws.Range("E24").Value = 28.89999961853028
$ws.Range("F24").Value = 35.75
$ws.Range("G24").Value = 23.5
$ws.Range("H24").Value = 9846715
$ws.Range("I24").Value = "KSCP"

$ws.Range("D25").Value = 25.5
$ws.Range("E25").Value = 22
$ws.Range("F25").Value = 31.5
$ws.Range("G25").Value = 21
$ws.Range("H25").Value = 9846715
$ws.Range("I25").Value = "KSCP"

$ws.Range("D26").Value = 16.20000076293945
$ws.Range("E26").Value = 12
$ws.Range("F26").Value = 16.25
$ws.Range("G26").Value = 10.64999961853027
$ws.Range("H26").Value = 9846715
$ws.Range("I26").Value = "KSCP"

$ws.Range("D27").Value = 7.570000171661377
$ws.Range("E27").Value = 18.38999938964844
$ws.Range("F27").Value = 20.45999908447266
$ws.Range("G27").Value = 5.409999847412109
$ws.Range("H27").Value = 9846715
$ws.Range("I27").Value = "KSCP"

$ws.Range("D28").Value = 13.46000003814697
$ws.Range("E28").Value = 10.51000022888184
$ws.Range("F28").Value = 17.78800010681152
$ws.Range("G28").Value = 10.03999996185303
$ws.Range("H28").Value = 9846715
$ws.Range("I28").Value = "KSCP"

$ws.Range("D29").Value = 2.75
$ws.Range("E29").Value = 4.989999771118164
$ws.Range("F29").Value = 5.699999809265137
$ws.Range("G29").Value = 2.450000047683716
$ws.Range("H29").Value = 9846715
$ws.Range("I29").Value = "KSCP"

$ws.Range("D30").Value = 4.860000133514404
$ws.Range("E30").Value = 5.949999809265137
$ws.Range("F30").Value = 10.14000034332275
$ws.Range("G30").Value = 4.735000133514404
$ws.Range("H30").Value = 9846715
$ws.Range("I30").Value = "KSCP"
